$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-03 Friday" "2025-01-04 Saturday"

Replace-Text "453×3=" "677×2="
Replace-Text "441×4=" "178×2="
Replace-Text "798×2=" "212×7="
Replace-Text "679×3=" "460×9="
Replace-Text "546×4=" "699×6="

Replace-Text "625×7=" "458×6="
Replace-Text "533×5=" "707×5="
Replace-Text "934×5=" "850×4="
Replace-Text "514×7=" "690×5="
Replace-Text "192×5=" "930×9="

Replace-Text "114×3=" "635×2="
Replace-Text "790×5=" "793×3="
Replace-Text "998×4=" "535×2="
Replace-Text "491×8=" "337×6="
Replace-Text "656×6=" "764×7="

Replace-Text "847×6=" "655×7="
Replace-Text "734×3=" "426×2="
Replace-Text "524×6=" "537×8="
Replace-Text "263×9=" "319×4="
Replace-Text "928×6=" "546×6="

Replace-Text "742×3=" "512×9="
Replace-Text "265×5=" "836×4="
Replace-Text "284×6=" "669×5="
Replace-Text "612×7=" "813×8="
Replace-Text "525×5=" "313×4="
